$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Update the "ModelParameterSheets" value for "TestScenario2" (row 3, column E)
# to include the newly added parameter sheet "Sheet, with comma".
$ws.Range("E3").Value = '"Global", "Aciclovir", "Sheet, with comma"'
